# Add a "Save" column in column H of the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - styled like the other header cells (bold, bordered),
# matching the style used by B1:G1. Copy the formatting from G1 so the
# existing cell style is reused instead of creating a new one.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Save values per row (2-42), corresponding to whether the "sum" (column G)
# crosses the save threshold.
$saveValues = @(0,0,0,1,1,0,1,0,0,1,0,0,0,0,0,0,0,0,1,1,0,0,0,1,0,0,0,0,0,1,0,0,0,1,1,1,0,0,0,0,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
